# Fix latency units in report sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header O2 label
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to Read Latency columns (I, J, K) for data rows 3-15
$cols = @("I", "J", "K")
for ($row = 3; $row -le 15; $row++) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value2
        $cell.Value = "$current msec"
    }
}
